# Applies the template_exp_en.docx edits described by the commit:
# "Verbesserungen: Scheindruck klappt, Pruefungsdruck klappt,
#  Einstiegsseite etwas erweitert."

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Nudge the anchored letterhead picture a little to the left
#    (wp:posOffset 1885315 -> 1884680 EMU == 148.45pt -> 148.4pt)
# ---------------------------------------------------------------
$d.Shapes(1).Left = 148.4

# ---------------------------------------------------------------
# 2. Slightly rebalance the results-table column widths
#    (column 1: 6239->6238 dxa, column 3: 786->785 dxa,
#     column 4: 561->562 dxa; column 2 stays at 1470 dxa)
# ---------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.Columns(1).Width = 311.9
$tbl.Columns(3).Width = 39.25
$tbl.Columns(4).Width = 28.1

# ---------------------------------------------------------------
# 3. Refresh the cached "today" date next to the signature line
# ---------------------------------------------------------------
$d.Content.Find.Execute("08.06.16", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "03.07.16", 2)

# ---------------------------------------------------------------
# 4. Swap out the academic advisor's name
# ---------------------------------------------------------------
$d.Content.Find.Execute("Jakubzik", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Normalverbraucher", 2)

# ---------------------------------------------------------------
# 5. Flip the "overflowPunct" compatibility flag on three styles
#    (ParagraphFormat.HangingPunctuation is this engine's surface
#    for the w:overflowPunct paragraph-property toggle)
# ---------------------------------------------------------------
$d.Styles("Normal").ParagraphFormat.HangingPunctuation = $false
$d.Styles("HDFuabsender").ParagraphFormat.HangingPunctuation = $true
$d.Styles("HDAufzhlung").ParagraphFormat.HangingPunctuation = $false

# ---------------------------------------------------------------
# 6. Add the two new (unused) list-label character styles that
#    follow ListLabel85 in the stylesheet
# ---------------------------------------------------------------
$ll86 = $d.Styles.Add("ListLabel86", 2)
$ll86.NameLocal = "ListLabel 86"
$ll86.Font.NameBi = "Symbol"

$ll87 = $d.Styles.Add("ListLabel87", 2)
$ll87.NameLocal = "ListLabel 87"
$ll87.Font.NameBi = "OpenSymbol"

Write-Output "edits applied"
